$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.113.00'
$ws.Range("E2").Value = '  -0.61%  '

$ws.Range("D3").Value = '3.007.77'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.57'
$ws.Range("E5").Value = '  +1.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.83'
$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '3.006.80'
$ws.Range("E8").Value = '  +0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.24'
$ws.Range("E10").Value = '  +7.24%  '

$ws.Range("E11").Value = '  +0.07%  '

$ws.Range("E12").Value = '  -0.65%  '

$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.33'
$ws.Range("E14").Value = '  -0.95%  '

$ws.Range("E15").Value = '  +2.53%  '

$ws.Range("D16").Value = '3.502.78'
$ws.Range("E16").Value = '  +0.01%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '62.099.63'
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.97'
$ws.Range("E18").Value = '  -1.54%  '

$ws.Range("D19").Value = '3.005.08'
$ws.Range("E19").Value = '  -0.11%  '

$ws.Range("E20").Value = '  -2.69%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.11'
$ws.Range("E21").Value = '  +1.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  -0.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.38'
$ws.Range("E23").Value = '  -0.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.14'
$ws.Range("E24").Value = '  +0.78%  '

$ws.Range("E25").Value = '  +10.69%  '

$ws.Range("E26").Value = '  +0.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.10'
$ws.Range("E27").Value = '  -1.98%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("E29").Value = '  +3.12%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.24'
$ws.Range("E31").Value = '  +3.07%  '

$ws.Range("E32").Value = '  -0.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.36'
$ws.Range("E33").Value = '  -2.64%  '

$ws.Range("E34").Value = '  +1.12%  '

$ws.Range("E35").Value = '  +4.57%  '

$ws.Range("E36").Value = '  -0.41%  '

$ws.Range("E37").Value = '  +0.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.15'
$ws.Range("E38").Value = '  -0.38%  '

$ws.Range("E39").Value = '  -3.49%  '

$ws.Range("E40").Value = '  -1.58%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.95'
$ws.Range("E41").Value = '  +1.40%  '

$ws.Range("E42").Value = '  +2.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '41.54'
$ws.Range("E43").Value = '  +10.65%  '

$ws.Range("E44").Value = '  +4.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '391.06'
$ws.Range("E45").Value = '  +0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0351'
$ws.Range("E46").Value = '  -1.94%  '

$ws.Range("D47").Value = '2.718.24'
$ws.Range("E47").Value = '  -0.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.67'
$ws.Range("E48").Value = '  +3.34%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.17'
$ws.Range("E50").Value = '  -1.23%  '

$ws.Range("E51").Value = '  -1.56%  '
